# Martin County Football Roster - "Add files via upload"
#
# The underlying edit removes the roster row for "Damian Cheek"
# (Number 18, Jersey 25, RB/LB, Sr.) which previously occupied row 13.
# In the saved workbook this shows up as the row-13 cells being cleared
# out entirely (so row 13 disappears from sheetData while every other
# row keeps its original row number), the "Damian Cheek" shared string
# being dropped, and the selection left on the now-empty A13:E13 range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Updated_Football_Roster")

# Select the row that is about to be removed (mirrors selecting the row
# in the UI before deleting its contents) and clear it out.
$row13 = $ws.Range("A13:E13")
$row13.Select()
$row13.ClearContents()
